$wb = $excel.ActiveWorkbook

# --- Sheet "Nombre Pelicula": add the two new movie-name rows ---
$wsNombre = $wb.Worksheets.Item("Nombre Pelicula")
$wsNombre.Range("A13").Value = "Testigo de cargoz"
$wsNombre.Range("A14").Value = "Testigo de cargo"

# --- Sheet "Genero": add the two new genre rows ---
$wsGenero = $wb.Worksheets.Item("Genero")
$wsGenero.Range("A6").Value = "Mafia"
$wsGenero.Range("A7").Value = "Documental"

# --- Sheet "Pelicula": fix the genre of an existing row and append a new one ---
$wsPelicula = $wb.Worksheets.Item("Pelicula")

# Correct "Testigo de cargo" genre from "Intriga" to "Documental"
$wsPelicula.Range("B2").Value = "Documental"

# Append "El Padrino. Parte 2" as a new row
$wsPelicula.Range("A12").Value = "El Padrino. Parte 2"
$wsPelicula.Range("B12").Value = "Mafia"
$wsPelicula.Range("C12").Value = "Francis Ford Coppola"
$wsPelicula.Range("D12").Value = "Estados Unidos"
$wsPelicula.Range("E12").Value = "Francis Ford Coppola"

# F12/G12 ("1974" / "8.9") look numeric, but need to be stored as text (like the
# rest of column F/G in this sheet). Stage them through a scratch cell formatted
# as text, then copy the already-typed value across so the destination keeps the
# text representation instead of being reinterpreted as a number.
$scratch = $wsPelicula.Range("Z1")
$scratch.NumberFormat = "@"

$scratch.Value = "1974"
$scratch.Copy($wsPelicula.Range("F12"))
$wsPelicula.Range("F12").ClearFormats()

$scratch.Value = "8.9"
$scratch.Copy($wsPelicula.Range("G12"))
$wsPelicula.Range("G12").ClearFormats()

$scratch.Clear()
